$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B23").Value = "Add application icon"
$ws.Range("B24").Value = "Delete downloaded file after download canceled"

$ws.Range("B25").Select()
